$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 162.5
$ws.Range("I33").Value = 181.33333
$ws.Range("J33").Value = 106
$ws.Range("K33").Value = 181.33333
$ws.Range("L33").Value = 106
$ws.Range("M33").Value = 47.66667000000001
$ws.Range("N33").Value = -564
$ws.Range("H58").Value = 435.25
$ws.Range("I58").Value = 435.25
$ws.Range("K58").Value = 1305.75
$ws.Range("M58").Value = -1155.75
$ws.Range("H82").Value = 778.6
$ws.Range("I82").Value = 778.6
$ws.Range("K82").Value = 2335.8
$ws.Range("M82").Value = -1929.8
$ws.Range("H85").Value = 778.6
$ws.Range("I85").Value = 778.6
$ws.Range("K85").Value = 2335.8
$ws.Range("M85").Value = -931.8000000000002
$ws.Range("H92").Value = 1885.8889
$ws.Range("I92").Value = 2029
$ws.Range("J92").Value = 1599.6666
$ws.Range("K92").Value = 2029
$ws.Range("L92").Value = 1599.6666
$ws.Range("M92").Value = -781
$ws.Range("N92").Value = -4095.6666
$ws.Range("H97").Value = 2276.2
$ws.Range("J97").Value = 2276.2
$ws.Range("L97").Value = 6828.599999999999
$ws.Range("N97").Value = -7820.599999999999
$ws.Range("H98").Value = 4686.8184
$ws.Range("I98").Value = 4270.6313
$ws.Range("J98").Value = 7322.6665
$ws.Range("K98").Value = 4270.6313
$ws.Range("L98").Value = 7322.6665
$ws.Range("M98").Value = -2772.6313
$ws.Range("N98").Value = -10318.6665
$ws.Range("H99").Value = 295.85715
$ws.Range("J99").Value = 200
$ws.Range("L99").Value = 600
$ws.Range("N99").Value = -3596
$ws.Range("H103").Value = 5103510.5
$ws.Range("I103").Value = 2263.5
$ws.Range("J103").Value = 11905174
$ws.Range("K103").Value = 6790.5
$ws.Range("L103").Value = 35715522
$ws.Range("M103").Value = -6204.5
$ws.Range("N103").Value = -35716694
$ws.Range("H104").Value = 799
$ws.Range("H113").Value = 4826.5
$ws.Range("J113").Value = 5377.75
$ws.Range("L113").Value = 5377.75
$ws.Range("N113").Value = -11885.75
$ws.Range("H118").Value = 498
$ws.Range("I118").Value = 497.5
$ws.Range("J118").Value = 500
$ws.Range("K118").Value = 1492.5
$ws.Range("L118").Value = 1500
$ws.Range("M118").Value = 164.5
$ws.Range("N118").Value = -4814
$ws.Range("H122").Value = 4686.8184
$ws.Range("I122").Value = 4270.6313
$ws.Range("J122").Value = 7322.6665
$ws.Range("K122").Value = 12811.8939
$ws.Range("L122").Value = 21967.9995
$ws.Range("M122").Value = -10361.8939
$ws.Range("N122").Value = -26867.9995
$ws.Range("H127").Value = 495.44446
$ws.Range("I127").Value = 310
$ws.Range("K127").Value = 930
$ws.Range("M127").Value = 4030
$ws.Range("H131").Value = 835771.5600000001
$ws.Range("I131").Value = 1112873.5
$ws.Range("J131").Value = 4466
$ws.Range("K131").Value = 3338620.5
$ws.Range("L131").Value = 13398
$ws.Range("M131").Value = -3333580.5
$ws.Range("N131").Value = -23478
$ws.Range("H132").Value = 6129.5415
$ws.Range("I132").Value = 6352.5654
$ws.Range("K132").Value = 19057.6962
$ws.Range("M132").Value = -16527.6962
$ws.Range("H137").Value = 1520817.5
$ws.Range("I137").Value = 2085092.2
$ws.Range("K137").Value = 6255276.6
$ws.Range("M137").Value = -6252726.6
$ws.Range("H138").Value = 2616.935
$ws.Range("I138").Value = 1234.375
$ws.Range("J138").Value = 2777.232
$ws.Range("K138").Value = 3703.125
$ws.Range("L138").Value = 8331.696
$ws.Range("M138").Value = 1436.875
$ws.Range("N138").Value = -18611.696
$ws.Range("H141").Value = 3197.9092
$ws.Range("I141").Value = 3197.9092
$ws.Range("K141").Value = 9593.7276
$ws.Range("M141").Value = -4413.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2376.9302
$ws.Range("I32").Value = 2362.1428
$ws.Range("K32").Value = 2362.1428
$ws.Range("M32").Value = -2075.1428
$ws.Range("H61").Value = 3723.7273
$ws.Range("I61").Value = 2098.5
$ws.Range("K61").Value = 2098.5
$ws.Range("M61").Value = -1886.5
$ws.Range("H122").Value = 2333.8215
$ws.Range("I122").Value = 2235.12
$ws.Range("J122").Value = 3156.3333
$ws.Range("K122").Value = 6705.36
$ws.Range("L122").Value = 9468.999899999999
$ws.Range("M122").Value = -4255.36
$ws.Range("N122").Value = -14368.9999
$ws.Range("H132").Value = 1651.5758
$ws.Range("I132").Value = 1049.2174
$ws.Range("K132").Value = 3147.6522
$ws.Range("M132").Value = -617.6522
$ws.Range("H136").Value = 3723.7273
$ws.Range("I136").Value = 2098.5
$ws.Range("K136").Value = 6295.5
$ws.Range("M136").Value = -3745.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 38000
$ws.Range("J64").Value = 38000
$ws.Range("L64").Value = 38000
$ws.Range("N64").Value = -38496
$ws.Range("H67").Value = 38000
$ws.Range("J67").Value = 38000
$ws.Range("L67").Value = 38000
$ws.Range("N67").Value = -39716
$ws.Range("H68").Value = 82295
$ws.Range("J68").Value = 82295
$ws.Range("L68").Value = 82295
$ws.Range("N68").Value = -83793
$ws.Range("H71").Value = 82295
$ws.Range("J71").Value = 82295
$ws.Range("L71").Value = 246885
$ws.Range("N71").Value = -254373
$ws.Range("H132").Value = 11909501
$ws.Range("I132").Value = 4424
$ws.Range("J132").Value = 20838308
$ws.Range("K132").Value = 13272
$ws.Range("L132").Value = 62514924
$ws.Range("M132").Value = -10742
$ws.Range("N132").Value = -62519984
$ws.Range("H134").Value = 4558.5264
$ws.Range("I134").Value = 5021.2
$ws.Range("K134").Value = 15063.6
$ws.Range("M134").Value = -12528.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 64103.383
$ws.Range("I9").Value = 82734
$ws.Range("J9").Value = 2001.3334
$ws.Range("K9").Value = 248202
$ws.Range("L9").Value = 6004.0002
$ws.Range("M9").Value = -247978
$ws.Range("N9").Value = -6452.0002
$ws.Range("H51").Value = 3500
$ws.Range("I51").Value = 1500
$ws.Range("K51").Value = 4500
$ws.Range("M51").Value = -4040
$ws.Range("H92").Value = 477.22223
$ws.Range("I92").Value = 265
$ws.Range("K92").Value = 795
$ws.Range("M92").Value = 453
$ws.Range("H109").Value = 3995.1428
$ws.Range("H132").Value = 1491.1666
$ws.Range("J132").Value = 1999.6666
$ws.Range("L132").Value = 17996.9994
$ws.Range("N132").Value = -23056.9994

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13412.333
$ws.Range("I40").Value = 15594.8
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 15594.8
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -15458.8
$ws.Range("N40").Value = -2772
$ws.Range("H68").Value = 2999.5
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 3999
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 3999
$ws.Range("N68").Value = -5497
$ws.Range("M68").Value = -1251
$ws.Range("H71").Value = 2999.5
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 3999
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 19995
$ws.Range("N71").Value = -27483
$ws.Range("M71").Value = -6256
$ws.Range("H136").Value = 4958.1113
$ws.Range("I136").Value = 3577.1
$ws.Range("J136").Value = 6684.375
$ws.Range("K136").Value = 10731.3
$ws.Range("L136").Value = 20053.125
$ws.Range("M136").Value = -8181.299999999999
$ws.Range("N136").Value = -25153.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7500
$ws.Range("I62").Value = 7500
$ws.Range("K62").Value = 7500
$ws.Range("M62").Value = -6876
$ws.Range("H65").Value = 7500
$ws.Range("I65").Value = 7500
$ws.Range("K65").Value = 37500
$ws.Range("M65").Value = -34380
$ws.Range("H126").Value = 2891.1428
$ws.Range("I126").Value = 2937
$ws.Range("J126").Value = 2830
$ws.Range("K126").Value = 8811
$ws.Range("L126").Value = 8490
$ws.Range("M126").Value = -6341
$ws.Range("N126").Value = -13430
$ws.Range("H132").Value = 1562.1794
$ws.Range("I132").Value = 1476
$ws.Range("J132").Value = 2596.3333
$ws.Range("K132").Value = 4428
$ws.Range("L132").Value = 7788.999899999999
$ws.Range("M132").Value = -1898
$ws.Range("N132").Value = -12848.9999
